# Apply updated "dSF" (column F) values for a set of rows in the
# active worksheet, per the repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = 3
    12 = 3
    15 = -2
    21 = 0
    24 = 0
    25 = -1
    29 = 3
    30 = -1
    34 = 4
    37 = -2
    43 = 1
    44 = 8
    46 = -5
    47 = -3
    48 = -5
    53 = -2
    54 = -6
    56 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
